# Apply odds updates to Sheet1, rows 2 and 6, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 6.5
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15
$ws.Range("AH2").Value = 19
$ws.Range("AJ2").Value = 21
$ws.Range("AO2").Value = 7
$ws.Range("AW2").Value = 8
$ws.Range("AZ2").Value = 126

# Row 6 updates
$ws.Range("G6").Value = 3.05
$ws.Range("I6").Value = 2.57
$ws.Range("J6").Value = 3.7
$ws.Range("K6").Value = 1.85
$ws.Range("L6").Value = 3.25
$ws.Range("Q6").Value = 2.45
$ws.Range("R6").Value = 1.42
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.18
$ws.Range("W6").Value = 7
$ws.Range("Z6").Value = 40
$ws.Range("AF6").Value = 110
$ws.Range("AK6").Value = 30
$ws.Range("AL6").Value = 26
$ws.Range("AM6").Value = 45
$ws.Range("AN6").Value = 4.7
$ws.Range("AO6").Value = 18
$ws.Range("AP6").Value = 28
$ws.Range("AR6").Value = 150
$ws.Range("AV6").Value = 80
$ws.Range("AW6").Value = 4.25
$ws.Range("AX6").Value = 14.5
$ws.Range("AY6").Value = 25
$ws.Range("AZ6").Value = 70
$ws.Range("BA6").Value = 120
$ws.Range("BB6").Value = 400
